$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4821145882335259
$ws.Range("H2").Value = 0.992

# Row 3
$ws.Range("D3").Value = 0.321570631223101
$ws.Range("G3").Value = 0.4821145882335259
$ws.Range("H3").Value = 0.992

# Row 4
$ws.Range("B4").Value = 0.2192987281846224
$ws.Range("D4").Value = 0.3777406617731509
$ws.Range("G4").Value = 0.4821145882335259
$ws.Range("H4").Value = 0.992

# Row 5
$ws.Range("G5").Value = 0.4821145882335259
$ws.Range("H5").Value = 0.992
